$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (45206 = 2023-10-07) for every
# data row (2..163). This update bumps that value to 45208 (2023-10-09) for all
# of them, leaving every other cell untouched.
$ws.Range("C2:C163").Value = 45208
